$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-08 06:46:07"

# --- Row 3: updated to former row-5 listing (5487838) ---
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【時給2,000円〜】フルスタックエンジニア募集|個人向けWebサービスの開発・保守・運用"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5487838"
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = "◆開発"

# --- Row 4: brand-new listing (5487908) ---
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "bubbleで構築したサイトの修正対応"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5487908"
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = "◇サイト"

# --- Row 5: updated to former row-6 listing (5487828), no skill summary ---
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "BigQuery+Looker Studioによる不動産マーケ分析ダッシュボード構築"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5487828"
$ws.Range("G5").Value = 25
$ws.Range("H5").ClearContents()

# --- Remove the now-obsolete trailing row 6 entirely ---
$ws.Rows.Item(6).Delete()

# --- Rebuild the hyperlinks so relationship targets match the new URLs ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5487791")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5487838")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5487908")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5487828")
